# Updates cryptos list figures (price/volume) and restores two swapped-row
# name/link pairs, per the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value. NumberFormat is forced to "@"
# (Text) before each write so numeric-looking strings (e.g. "1.000",
# "27.278.08") are stored verbatim as text instead of being parsed as
# numbers, matching the workbook's original inline-string cell type.
$updates = [ordered]@{
    "D2" = "27.278.08"
    "E2" = "  -1.01%  "
    "D3" = "1.785.46"
    "E3" = "  -2.47%  "
    "D4" = "1.004"
    "E4" = "  +0.12%  "
    "D5" = "339.85"
    "E5" = "  -1.37%  "
    "D6" = "1.000"
    "E6" = "  +0.07%  "
    "D7" = "0.3971"
    "E7" = "  +3.51%  "
    "D8" = "0.3438"
    "E8" = "  -3.27%  "
    "D9" = "48.03"
    "E9" = "  -4.40%  "
    "D10" = "1.190"
    "E10" = "  -4.45%  "
    "D11" = "0.07434"
    "E11" = "  -4.44%  "
    "D12" = "1.002"
    "E12" = "  +0.07%  "
    "E13" = "  -2.70%  "
    "D14" = "6.451"
    "E14" = "  -2.71%  "
    "D15" = "1.785.12"
    "E15" = "  -2.51%  "
    "D16" = "7.097"
    "E16" = "  -2.14%  "
    "D17" = "0.00001089"
    "E17" = "  -3.60%  "
    "D18" = "0.06672"
    "E18" = "  -1.21%  "
    "D19" = "83.97"
    "E19" = "  -3.55%  "
    "D20" = "1.000"
    "E20" = "  -0.01%  "
    "D21" = "17.68"
    "E21" = "  +0.20%  "
    "D22" = "6.487"
    "E22" = "  -1.32%  "
    "D23" = "27.273.14"
    "E23" = "  -1.04%  "
    "D24" = "12.31"
    "E24" = "  -6.84%  "
    "D25" = "2.377"
    "E25" = "  -3.75%  "
    "D26" = "1.471"
    "E26" = "  -2.12%  "
    "B27" = "EthereumClassic"
    "C27" = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
    "D27" = "21.07"
    "E27" = "  -4.91%  "
    "B28" = "LidoDAOToken"
    "C28" = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
    "D28" = "2.487"
    "E28" = "  -9.04%  "
    "D29" = "156.58"
    "E29" = "  +1.83%  "
    "D30" = "1.985.44"
    "E30" = "  -2.53%  "
    "D31" = "134.63"
    "E31" = "  -0.86%  "
    "D32" = "4.041"
    "E32" = "  -1.28%  "
    "D33" = "5.979"
    "E33" = "  -6.48%  "
    "D34" = "0.08774"
    "E34" = "  -0.58%  "
    "D35" = "12.97"
    "E35" = "  -7.42%  "
    "D36" = "1.621"
    "E36" = "  -4.40%  "
    "D37" = "5.384"
    "E37" = "  -4.97%  "
    "B38" = "VeChain"
    "C38" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "D38" = "0.02389"
    "E38" = "  -0.92%  "
    "B39" = "TheSandbox"
    "C39" = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
    "D39" = "0.6802"
    "E39" = "  -4.08%  "
    "D40" = "0.06399"
    "E40" = "  -2.19%  "
    "D41" = "0.2195"
    "E41" = "  -3.09%  "
    "D42" = "1.248"
    "E42" = "  -5.52%  "
    "D43" = "8.406"
    "E43" = "  -8.07%  "
    "D44" = "14.30"
    "E44" = "  -3.21%  "
    "D45" = "1.0000"
    "E45" = "  +0.04%  "
    "D46" = "0.6374"
    "E46" = "  -4.20%  "
    "D47" = "3.873"
    "E47" = "  -2.27%  "
    "D48" = "2.127"
    "E48" = "  -3.37%  "
    "D49" = "131.92"
    "E49" = "  -1.19%  "
    "D50" = "0.07127"
    "E50" = "  -2.74%  "
    "D51" = "78.64"
    "E51" = "  -3.42%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
